# Graafikuseedija - anonymize the respondent data:
#  - clear the "Email" formula column (D) for every data row
#  - replace the "Name" column (E) values with a redacted placeholder
#  - move the active cell selection to D27
#  - nudge the workbook window geometry (best effort)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 25

# Clear the computed Email column (was a formula referencing the Name column)
$ws.Range("D2:D" + $lastRow).ClearContents()

# Redact every respondent's name
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = "?????"
}

# Match the author's final selection / window layout
[void]$ws.Range("D27").Select()

try {
    $win = $excel.ActiveWindow
    $win.Left = 3180
    $win.Top = 300
    $win.Width = 23205
    $win.Height = 14550
} catch {
}
